# Add files via upload
# - Fill the newly added "Snatch" / "Clean and Jerk" (C:D) columns with 0
#   on the SFM and SFF sheets.
# - Update the remembered selection on each sheet.
# - Make "SFF" the active sheet/tab (was "ScoreF").

$wb = $excel.ActiveWorkbook

# SFM sheet: add zeroed C2:D7 and leave that range selected.
$sfm = $wb.Worksheets.Item("SFM")
$sfm.Range("C2:D7").Value = 0
[void]$sfm.Range("C2:D7").Select()

# SFF sheet: add zeroed C2:D7, select D4, and make it the active tab.
$sff = $wb.Worksheets.Item("SFF")
$sff.Range("C2:D7").Value = 0
[void]$sff.Range("D4").Select()
$sff.Activate()
